# Répartition.xlsx update
# - Mark "Steve" as also responsible for "Fusion Personnages/Ennemis" (row 6)
#   and "Ennemis (State/Observer)" (row 7)
# - Add two new tasks at the bottom of the table:
#     "Factory Characters" -> Steve   (row 11)
#     "Factory Surfaces"   -> Margaux (row 12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "Steve"
$ws.Range("C7").Value = "Steve"

$ws.Range("B11").Value = "Factory Characters"
$ws.Range("C11").Value = "Steve"

$ws.Range("B12").Value = "Factory Surfaces"
$ws.Range("C12").Value = "Margaux"

# Move the active selection down below the table, matching the saved view.
[void]$ws.Range("C13").Select()
